$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 15
$ws.Cells.Item(14, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4122)

$ws.Cells.Item($row, 1).Value = 42620.888136574074
$ws.Cells.Item($row, 2).Value = 16
$ws.Cells.Item($row, 3).Value = 62
$ws.Cells.Item($row, 4).Value = 36
$ws.Cells.Item($row, 5).Value = 62
$ws.Cells.Item($row, 6).Value = 23
$ws.Cells.Item($row, 7).Value = 34074
$ws.Cells.Item($row, 8).Value = 27720
$ws.Cells.Item($row, 9).Value = 3249
$ws.Cells.Item($row, 10).Value = 393
$ws.Cells.Item($row, 11).Value = 225
$ws.Cells.Item($row, 12).Value = 30
$ws.Cells.Item($row, 13).Value = 9
$ws.Cells.Item($row, 14).Value = "Bag"
